$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- view: scroll / selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 157
$win.ScrollColumn = 1
[void]$ws.Range("B487").Select()

# --- template rows used to clone formatting from ---
# "fill" look (blue-ish highlighted row): styles 9/10/11 on A/B/C/D
$fillTemplate = $ws.Range("A4:D4")
# "note" look (yellow-ish "Inp" row): styles 21/21/21/22 on A/B/C/D
$noteTemplate = $ws.Range("A17:D17")

function Set-RowFormat($row, $template) {
    $template.Copy()
    $ws.Range("A" + $row + ":D" + $row).PasteSpecial(-4122)
}

function Set-Row($row, $style, $aVal, $bVal) {
    if ($style -eq "fill") {
        Set-RowFormat $row $fillTemplate
    } else {
        Set-RowFormat $row $noteTemplate
    }
    $ws.Range("A" + $row).Value = $aVal
    if ($bVal) {
        $ws.Range("B" + $row).Value = $bVal
    } else {
        $ws.Range("B" + $row).Value = ""
    }
}

Set-Row 56  "fill" "maths" "Done"
Set-Row 151 "fill" "TEdgeSegment" "Done"
Set-Row 164 "fill" "maths" "Done"
Set-Row 189 "fill" "maths" "Done"
Set-Row 213 "fill" "maths" "Done"

$ws.Range("B259").Value = "Inp"

Set-Row 278 "note" "proj" "Done"
Set-Row 279 "note" "proj" "Done"
Set-Row 280 "note" "proj" "Done"
Set-Row 281 "note" "proj" "Done"
Set-Row 282 "note" "proj" "Done"
Set-Row 283 "fill" "maths" "Done"
Set-Row 284 "fill" "maths" "Done"

Set-Row 412 "fill" "TCircle" $null
Set-Row 413 "fill" "TCircle" "Done"
Set-Row 414 "fill" "TCircle" $null
Set-Row 415 "fill" "TCircle" "Done"

Set-Row 416 "note" "TCollisionComponent" "Done"
Set-Row 417 "note" "TCollisionComponent" $null
Set-Row 418 "note" "TCollisionComponent" $null
Set-Row 419 "note" "TCollisionComponent" "Done"
Set-Row 420 "note" "TCollisionComponent" "Done"

Set-Row 454 "fill" "TEdgeSegment" "Done"

Set-Row 487 "fill" "timer" $null
Set-Row 488 "fill" "timer" "inp"
Set-Row 489 "fill" "timer" $null
Set-Row 490 "fill" "timer" $null
Set-Row 491 "fill" "timer" $null

Set-Row 532 "fill" "Tline" $null
Set-Row 533 "fill" "Tline" "Done"
Set-Row 534 "fill" "Tline" "Done"
Set-Row 535 "fill" "Tline" $null
Set-Row 536 "fill" "Tline" "Done"
Set-Row 537 "fill" "Tline" "Done"

$ws.Range("B558").Value = "Done"
